$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colA = @(6,12,18,24,30,36,42,48,54,60,66,72,78,84,90,96,102,108,114,120,126,132,138,144,150,156,162,168,174,180,186,192,198,204,210,216,222,228,234,240,246,252,258,264,270,276,282,288,294,300)
$colC = @(90.48,89.7,92.75,91.44,91.69,90.23,91.44,91.44,91.25,91.21,91.46,92.42,90.48,90.98,92.46,91.18,91.69,92.42,90.45,90.7,91.18,91.46,90.7,92.95,90.5,91.21,91.71,91.67,90.25,92.7,91.18,91.5,92.96,91.46,90.68,90.5,91.21,90.73,92.46,91.48,91.27,92,91.46,90.25,91.75,91.94,90.95,58.72,88.89,91.25)

for ($i = 0; $i -lt $colA.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $colA[$i]
    $ws.Cells.Item($row, 3).Value = $colC[$i]
}

[void]$ws.Range("E18").Select()
